$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1914.7273
$ws.Range("I15").Value = 1914.7273
$ws.Range("K15").Value = 5744.1819
$ws.Range("M15").Value = -5575.1819
$ws.Range("H28").Value = 605.2
$ws.Range("I28").Value = 584.7778
$ws.Range("K28").Value = 584.7778
$ws.Range("M28").Value = -99.77779999999996
$ws.Range("H29").Value = 2036
$ws.Range("I29").Value = 1381.6666
$ws.Range("J29").Value = 3999
$ws.Range("K29").Value = 4144.9998
$ws.Range("L29").Value = 11997
$ws.Range("M29").Value = -3863.9998
$ws.Range("N29").Value = -12559
$ws.Range("H40").Value = 2124.75
$ws.Range("I40").Value = 2166.3333
$ws.Range("K40").Value = 2166.3333
$ws.Range("M40").Value = -1991.3333
$ws.Range("H43").Value = 4266.6665
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 4266.6665
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 4266.6665
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -4404.6665
$ws.Range("H51").Value = 12250.75
$ws.Range("J51").Value = 11501
$ws.Range("L51").Value = 11501
$ws.Range("N51").Value = -12469
$ws.Range("H61").Value = 1050
$ws.Range("I61").Value = 1050
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3150
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2978
$ws.Range("N61").Value = ""
$ws.Range("H70").Value = 2328.1428
$ws.Range("I70").Value = 2074
$ws.Range("J70").Value = 2667
$ws.Range("K70").Value = 6222
$ws.Range("L70").Value = 8001
$ws.Range("M70").Value = -5952
$ws.Range("N70").Value = -8541
$ws.Range("H73").Value = 2328.1428
$ws.Range("I73").Value = 2074
$ws.Range("J73").Value = 2667
$ws.Range("K73").Value = 6222
$ws.Range("L73").Value = 8001
$ws.Range("M73").Value = -5286
$ws.Range("N73").Value = -9873
$ws.Range("H96").Value = 2025.5264
$ws.Range("I96").Value = 1208.2307
$ws.Range("J96").Value = 3796.3333
$ws.Range("K96").Value = 3624.6921
$ws.Range("L96").Value = 11388.9999
$ws.Range("M96").Value = -2251.6921
$ws.Range("N96").Value = -14134.9999
$ws.Range("H116").Value = 5485.4165
$ws.Range("I116").Value = 5443.9414
$ws.Range("J116").Value = 5586.143
$ws.Range("K116").Value = 5443.9414
$ws.Range("L116").Value = 5586.143
$ws.Range("M116").Value = -2001.9414
$ws.Range("N116").Value = -12470.143
$ws.Range("H129").Value = 1089.9642
$ws.Range("I129").Value = 1967.1428
$ws.Range("J129").Value = 797.5714
$ws.Range("K129").Value = 5901.428400000001
$ws.Range("L129").Value = 2392.7142
$ws.Range("M129").Value = -901.4284000000007
$ws.Range("N129").Value = -12392.7142
$ws.Range("H132").Value = 2873.3691
$ws.Range("I132").Value = 2342.1694
$ws.Range("K132").Value = 7026.5082
$ws.Range("M132").Value = -4496.5082
$ws.Range("H135").Value = 1480.2195
$ws.Range("I135").Value = 993.9697
$ws.Range("J135").Value = 3486
$ws.Range("K135").Value = 8945.7273
$ws.Range("L135").Value = 31374
$ws.Range("M135").Value = -6410.7273
$ws.Range("N135").Value = -36444
$ws.Range("H137").Value = 4905.1665
$ws.Range("I137").Value = 1257.7778
$ws.Range("J137").Value = 8552.556
$ws.Range("K137").Value = 3773.3334
$ws.Range("L137").Value = 25657.668
$ws.Range("M137").Value = -1223.3334
$ws.Range("N137").Value = -30757.668
$ws.Range("H141").Value = 14766.333
$ws.Range("I141").Value = 14766.333
$ws.Range("K141").Value = 44298.999
$ws.Range("M141").Value = -39118.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 849.9737
$ws.Range("I2").Value = 739.7
$ws.Range("J2").Value = 1263.5
$ws.Range("K2").Value = 739.7
$ws.Range("L2").Value = 1263.5
$ws.Range("M2").Value = -626.7
$ws.Range("N2").Value = -1489.5
$ws.Range("H6").Value = 5000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5346
$ws.Range("H31").Value = 8356.700000000001
$ws.Range("I31").Value = 7063
$ws.Range("K31").Value = 7063
$ws.Range("M31").Value = -6769
$ws.Range("H32").Value = 9659.960999999999
$ws.Range("I32").Value = 5891.625
$ws.Range("J32").Value = 23363
$ws.Range("K32").Value = 5891.625
$ws.Range("L32").Value = 23363
$ws.Range("M32").Value = -5604.625
$ws.Range("N32").Value = -23937
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = ""
$ws.Range("H45").Value = 2433
$ws.Range("I45").Value = 1133.3636
$ws.Range("K45").Value = 1133.3636
$ws.Range("M45").Value = -756.3635999999999
$ws.Range("H61").Value = 4272.2324
$ws.Range("I61").Value = 3852.7896
$ws.Range("K61").Value = 3852.7896
$ws.Range("M61").Value = -3640.7896
$ws.Range("H74").Value = 102184.5
$ws.Range("I74").Value = 125793.125
$ws.Range("K74").Value = 125793.125
$ws.Range("M74").Value = -124919.125
$ws.Range("H77").Value = 102184.5
$ws.Range("I77").Value = 125793.125
$ws.Range("K77").Value = 628965.625
$ws.Range("M77").Value = -624597.625
$ws.Range("H102").Value = 1781.7084
$ws.Range("J102").Value = 2988.25
$ws.Range("L102").Value = 2988.25
$ws.Range("N102").Value = -6232.25
$ws.Range("H116").Value = 849.9737
$ws.Range("I116").Value = 739.7
$ws.Range("J116").Value = 1263.5
$ws.Range("K116").Value = 739.7
$ws.Range("L116").Value = 1263.5
$ws.Range("M116").Value = 1554.3
$ws.Range("N116").Value = -5851.5
$ws.Range("H132").Value = 52919.3
$ws.Range("I132").Value = 55547.156
$ws.Range("J132").Value = 2990
$ws.Range("K132").Value = 166641.468
$ws.Range("L132").Value = 8970
$ws.Range("M132").Value = -164111.468
$ws.Range("N132").Value = -14030
$ws.Range("H136").Value = 4272.2324
$ws.Range("I136").Value = 3852.7896
$ws.Range("K136").Value = 11558.3688
$ws.Range("M136").Value = -9008.3688
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 849.9737
$ws.Range("I3").Value = 739.7
$ws.Range("J3").Value = 1263.5
$ws.Range("K3").Value = 739.7
$ws.Range("L3").Value = 1263.5
$ws.Range("M3").Value = -625.7
$ws.Range("N3").Value = -1491.5
$ws.Range("H6").Value = 39999
$ws.Range("J6").Value = 39999
$ws.Range("L6").Value = 39999
$ws.Range("N6").Value = -40225
$ws.Range("H20").Value = 1686.6666
$ws.Range("I20").Value = 1404.6154
$ws.Range("J20").Value = 2420
$ws.Range("K20").Value = 1404.6154
$ws.Range("L20").Value = 2420
$ws.Range("M20").Value = -1157.6154
$ws.Range("N20").Value = -2914
$ws.Range("H81").Value = 539434
$ws.Range("J81").Value = 539434
$ws.Range("L81").Value = 539434
$ws.Range("N81").Value = -541556
$ws.Range("H84").Value = 539434
$ws.Range("J84").Value = 539434
$ws.Range("L84").Value = 1618302
$ws.Range("N84").Value = -1628910
$ws.Range("H94").Value = 1819.6
$ws.Range("I94").Value = 1305.1818
$ws.Range("J94").Value = 3234.25
$ws.Range("K94").Value = 1305.1818
$ws.Range("L94").Value = 3234.25
$ws.Range("M94").Value = -854.1818000000001
$ws.Range("N94").Value = -4136.25
$ws.Range("H105").Value = 5047.3
$ws.Range("I105").Value = 4913.5
$ws.Range("J105").Value = 5248
$ws.Range("K105").Value = 4913.5
$ws.Range("L105").Value = 5248
$ws.Range("M105").Value = -3166.5
$ws.Range("N105").Value = -8742
$ws.Range("H116").Value = 74331
$ws.Range("J116").Value = 74331
$ws.Range("L116").Value = 74331
$ws.Range("N116").Value = -83509
$ws.Range("H134").Value = 2475.3555
$ws.Range("J134").Value = 3631.5
$ws.Range("L134").Value = 10894.5
$ws.Range("N134").Value = -15964.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3617.2727
$ws.Range("I31").Value = 3183.35
$ws.Range("K31").Value = 3183.35
$ws.Range("M31").Value = -2888.35
$ws.Range("H34").Value = 3617.2727
$ws.Range("I34").Value = 3183.35
$ws.Range("K34").Value = 3183.35
$ws.Range("M34").Value = -2981.35
$ws.Range("H55").Value = 31832.666
$ws.Range("I55").Value = 31832.666
$ws.Range("K55").Value = 31832.666
$ws.Range("M55").Value = -31517.666
$ws.Range("H58").Value = 36196.965
$ws.Range("I58").Value = 42800.44
$ws.Range("J58").Value = 3179.6
$ws.Range("K58").Value = 42800.44
$ws.Range("L58").Value = 3179.6
$ws.Range("M58").Value = -42597.44
$ws.Range("N58").Value = -3585.6
$ws.Range("H132").Value = 4553.875
$ws.Range("I132").Value = 4698.769
$ws.Range("J132").Value = 3926
$ws.Range("K132").Value = 14096.307
$ws.Range("L132").Value = 11778
$ws.Range("M132").Value = -11566.307
$ws.Range("N132").Value = -16838
$ws.Range("H134").Value = 43097.88
$ws.Range("I134").Value = 46566.39
$ws.Range("J134").Value = 3210
$ws.Range("K134").Value = 139699.17
$ws.Range("L134").Value = 9630
$ws.Range("M134").Value = -137164.17
$ws.Range("N134").Value = -14700
$ws.Range("H136").Value = 36196.965
$ws.Range("I136").Value = 42800.44
$ws.Range("J136").Value = 3179.6
$ws.Range("K136").Value = 128401.32
$ws.Range("L136").Value = 9538.799999999999
$ws.Range("M136").Value = -125851.32
$ws.Range("N136").Value = -14638.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1966
$ws.Range("I3").Value = 1966
$ws.Range("K3").Value = 5898
$ws.Range("M3").Value = -5786
$ws.Range("H12").Value = 6.5555553
$ws.Range("I12").Value = 2.75
$ws.Range("J12").Value = 9.6
$ws.Range("K12").Value = 8.25
$ws.Range("L12").Value = 28.8
$ws.Range("M12").Value = 164.75
$ws.Range("N12").Value = -374.8
$ws.Range("H55").Value = 5080
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5080
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 15240
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -15594
$ws.Range("H62").Value = 7132.6
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7132.6
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 21397.8
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -22769.8
$ws.Range("H64").Value = 3913.5454
$ws.Range("J64").Value = 4611
$ws.Range("L64").Value = 13833
$ws.Range("N64").Value = -14373
$ws.Range("H65").Value = 7132.6
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7132.6
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 64193.4
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -71057.39999999999
$ws.Range("H67").Value = 3913.5454
$ws.Range("J67").Value = 4611
$ws.Range("L67").Value = 13833
$ws.Range("N67").Value = -15705
$ws.Range("H68").Value = 712.6667
$ws.Range("I68").Value = 727.4545000000001
$ws.Range("K68").Value = 2182.3635
$ws.Range("M68").Value = -1371.3635
$ws.Range("H71").Value = 712.6667
$ws.Range("I71").Value = 727.4545000000001
$ws.Range("K71").Value = 6547.0905
$ws.Range("M71").Value = -2491.0905
$ws.Range("H80").Value = 101436.25
$ws.Range("I80").Value = 1750
$ws.Range("J80").Value = 134665
$ws.Range("K80").Value = 5250
$ws.Range("L80").Value = 403995
$ws.Range("M80").Value = -4314
$ws.Range("N80").Value = -405867
$ws.Range("H82").Value = 6749.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 6749.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 20248.5
$ws.Range("M82").Value = ""
$ws.Range("N82").Value = -21060.5
$ws.Range("H83").Value = 101436.25
$ws.Range("I83").Value = 1750
$ws.Range("J83").Value = 134665
$ws.Range("K83").Value = 15750
$ws.Range("L83").Value = 1211985
$ws.Range("M83").Value = -11070
$ws.Range("N83").Value = -1221345
$ws.Range("H85").Value = 6749.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 6749.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 20248.5
$ws.Range("M85").Value = ""
$ws.Range("N85").Value = -23056.5
$ws.Range("H87").Value = 9742.6
$ws.Range("I87").Value = 9742.6
$ws.Range("K87").Value = 29227.8
$ws.Range("M87").Value = -27979.8
$ws.Range("H90").Value = 9742.6
$ws.Range("I90").Value = 9742.6
$ws.Range("K90").Value = 87683.40000000001
$ws.Range("M90").Value = -81443.40000000001
$ws.Range("H122").Value = 480.45834
$ws.Range("J122").Value = 560.5
$ws.Range("L122").Value = 5044.5
$ws.Range("N122").Value = -9944.5
$ws.Range("H132").Value = 958.1429000000001
$ws.Range("I132").Value = 879.8889
$ws.Range("K132").Value = 7919.0001
$ws.Range("M132").Value = -5389.0001
$ws.Range("H133").Value = 8000
$ws.Range("J133").Value = 9400
$ws.Range("L133").Value = 28200
$ws.Range("N133").Value = -38320
$ws.Range("H138").Value = 1894
$ws.Range("I138").Value = 1547.7142
$ws.Range("K138").Value = 4643.142599999999
$ws.Range("M138").Value = 496.8574000000008
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 12399
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 12399
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 12399
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -12737
$ws.Range("H70").Value = 7811
$ws.Range("I70").Value = 7833
$ws.Range("K70").Value = 7833
$ws.Range("M70").Value = -7563
$ws.Range("H73").Value = 7811
$ws.Range("I73").Value = 7833
$ws.Range("K73").Value = 7833
$ws.Range("M73").Value = -6897
$ws.Range("H80").Value = 4577.6
$ws.Range("I80").Value = 4301
$ws.Range("K80").Value = 4301
$ws.Range("M80").Value = -3303
$ws.Range("H83").Value = 4577.6
$ws.Range("I83").Value = 4301
$ws.Range("K83").Value = 21505
$ws.Range("M83").Value = -16513
$ws.Range("H97").Value = 925.7
$ws.Range("I97").Value = 905.44446
$ws.Range("J97").Value = 1108
$ws.Range("K97").Value = 905.44446
$ws.Range("L97").Value = 1108
$ws.Range("M97").Value = -409.44446
$ws.Range("N97").Value = -2100
$ws.Range("H102").Value = 3657.818
$ws.Range("I102").Value = 2693
$ws.Range("K102").Value = 2693
$ws.Range("M102").Value = -1071
$ws.Range("H113").Value = 1557.5333
$ws.Range("I113").Value = 1331.1111
$ws.Range("J113").Value = 1897.1666
$ws.Range("K113").Value = 1331.1111
$ws.Range("L113").Value = 1897.1666
$ws.Range("M113").Value = 838.8888999999999
$ws.Range("N113").Value = -6237.1666
$ws.Range("H122").Value = 2526.6875
$ws.Range("I122").Value = 1706.1818
$ws.Range("J122").Value = 4331.8
$ws.Range("K122").Value = 5118.5454
$ws.Range("L122").Value = 12995.4
$ws.Range("M122").Value = -2668.5454
$ws.Range("N122").Value = -17895.4
$ws.Range("H132").Value = 34222.625
$ws.Range("I132").Value = 36759.484
$ws.Range("J132").Value = 9699.666999999999
$ws.Range("K132").Value = 110278.452
$ws.Range("L132").Value = 29099.001
$ws.Range("M132").Value = -107748.452
$ws.Range("N132").Value = -34159.001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9396.559999999999
$ws.Range("I7").Value = 11795.647
$ws.Range("J7").Value = 4298.5
$ws.Range("K7").Value = 11795.647
$ws.Range("L7").Value = 4298.5
$ws.Range("M7").Value = -11683.647
$ws.Range("N7").Value = -4522.5
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("H55").Value = 1622.0769
$ws.Range("I55").Value = 1526.6111
$ws.Range("J55").Value = 1836.875
$ws.Range("K55").Value = 1526.6111
$ws.Range("L55").Value = 1836.875
$ws.Range("M55").Value = -1353.6111
$ws.Range("N55").Value = -2182.875
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H93").Value = 1200.3889
$ws.Range("I93").Value = 1065.1177
$ws.Range("K93").Value = 1065.1177
$ws.Range("M93").Value = 182.8823
$ws.Range("H94").Value = 36000
$ws.Range("J94").Value = 36000
$ws.Range("L94").Value = 36000
$ws.Range("N94").Value = -37352
$ws.Range("H100").Value = 4274.4443
$ws.Range("I100").Value = 4185
$ws.Range("K100").Value = 4185
$ws.Range("M100").Value = -3644
$ws.Range("H122").Value = 3889.261
$ws.Range("I122").Value = 3281.125
$ws.Range("J122").Value = 4552.6816
$ws.Range("K122").Value = 9843.375
$ws.Range("L122").Value = 13658.0448
$ws.Range("M122").Value = -7393.375
$ws.Range("N122").Value = -18558.0448
$ws.Range("H126").Value = 9396.559999999999
$ws.Range("I126").Value = 11795.647
$ws.Range("J126").Value = 4298.5
$ws.Range("K126").Value = 35386.94100000001
$ws.Range("L126").Value = 12895.5
$ws.Range("M126").Value = -32916.94100000001
$ws.Range("N126").Value = -17835.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
$ws.Range("H136").Value = 3437
$ws.Range("I136").Value = 3437
$ws.Range("K136").Value = 10311
$ws.Range("M136").Value = -7761
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 49264.668
$ws.Range("I58").Value = 49264.668
$ws.Range("K58").Value = 49264.668
$ws.Range("M58").Value = -48956.668
$ws.Range("H96").Value = 4613.4287
$ws.Range("I96").Value = 3358.8
$ws.Range("K96").Value = 3358.8
$ws.Range("M96").Value = -1985.8
$ws.Range("H122").Value = 1373.1316
$ws.Range("I122").Value = 1222.2572
$ws.Range("J122").Value = 3133.3333
$ws.Range("K122").Value = 3666.7716
$ws.Range("L122").Value = 9399.999899999999
$ws.Range("M122").Value = -1216.7716
$ws.Range("N122").Value = -14299.9999
$ws.Range("H132").Value = 23761.912
$ws.Range("I132").Value = 25228.162
$ws.Range("K132").Value = 75684.486
$ws.Range("M132").Value = -73154.486
$ws.Range("H136").Value = 4412.467
$ws.Range("I136").Value = 4479.136
$ws.Range("K136").Value = 13437.408
$ws.Range("M136").Value = -10887.408
